$d = $word.ActiveDocument

function Add-FigureCaption($searchText, $bmName) {
    # Locate the image-caption paragraph's text (e.g. "testimg") which sits
    # at the very start of its paragraph.
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $startPos = $rng.Start

    # "Figure "
    $insPt = $d.Range($startPos, $startPos)
    $insPt.InsertBefore("Figure ")

    # SEQ field (renders as the figure number, auto-numbered by Word)
    $fldPos = $d.Range($startPos + 7, $startPos + 7)
    $fld = $d.Fields.Add($fldPos, -1, "SEQ  \* ARABIC", $false)

    # Re-find the original text to know exactly where the field result ends.
    $rng2 = $d.Content
    $found2 = $rng2.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $afterField = $rng2.Start

    # ":" then " "
    $c1 = $d.Range($afterField, $afterField)
    $c1.InsertBefore(":")
    $c2 = $d.Range($afterField + 1, $afterField + 1)
    $c2.InsertBefore(" ")

    # Bookmark wrapping "Figure <n>: "
    $bmRange = $d.Range($startPos, $afterField + 2)
    $d.Bookmarks.Add($bmName, $bmRange)
}

Add-FigureCaption "testimg" "fig1"
Add-FigureCaption "2testimg" "fig2"
Add-FigureCaption "3testimg" "fig3"
Add-FigureCaption "4testimg" "fig4"
